# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice / NQ / HQ prices, leve profit calcs) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with freshly fetched data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 897.3
$ws.Range("I39").Value = 211.6
$ws.Range("K39").Value = 634.8
$ws.Range("M39").Value = -338.8

$ws.Range("H64").Value = 8799
$ws.Range("J64").Value = 11444.444
$ws.Range("L64").Value = 11444.444
$ws.Range("N64").Value = -11940.444

$ws.Range("H67").Value = 8799
$ws.Range("J67").Value = 11444.444
$ws.Range("L67").Value = 11444.444
$ws.Range("N67").Value = -13160.444

$ws.Range("H80").Value = 343.33334
$ws.Range("I80").Value = 292.1
$ws.Range("K80").Value = 876.3000000000001
$ws.Range("M80").Value = 121.6999999999999

$ws.Range("H83").Value = 343.33334
$ws.Range("I83").Value = 292.1
$ws.Range("K83").Value = 2628.9
$ws.Range("M83").Value = 2363.1

$ws.Range("H96").Value = 544.5714
$ws.Range("I96").Value = 453.5
$ws.Range("J96").Value = 772.25
$ws.Range("K96").Value = 1360.5
$ws.Range("L96").Value = 2316.75
$ws.Range("M96").Value = 12.5
$ws.Range("N96").Value = -5062.75

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H125").Value = 17806.953
$ws.Range("I125").Value = 24666.445
$ws.Range("J125").Value = 12662.333
$ws.Range("K125").Value = 221998.005
$ws.Range("L125").Value = 113960.997
$ws.Range("M125").Value = -219538.005
$ws.Range("N125").Value = -118880.997

$ws.Range("H131").Value = 117106.11
$ws.Range("I131").Value = 128994.375
$ws.Range("J131").Value = 22000
$ws.Range("K131").Value = 386983.125
$ws.Range("L131").Value = 66000
$ws.Range("M131").Value = -381943.125
$ws.Range("N131").Value = -76080

$ws.Range("H132").Value = 1978.6818
$ws.Range("I132").Value = 1956.1052
$ws.Range("J132").Value = 2121.6667
$ws.Range("K132").Value = 5868.3156
$ws.Range("L132").Value = 6365.000100000001
$ws.Range("M132").Value = -3338.3156
$ws.Range("N132").Value = -11425.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2633
$ws.Range("I122").Value = 2633
$ws.Range("K122").Value = 7899
$ws.Range("M122").Value = -5449

$ws.Range("H132").Value = 3512.025
$ws.Range("I132").Value = 2907.3948
$ws.Range("K132").Value = 8722.1844
$ws.Range("M132").Value = -6192.1844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1242.826
$ws.Range("I86").Value = 1388.7368
$ws.Range("J86").Value = 549.75
$ws.Range("K86").Value = 1388.7368
$ws.Range("L86").Value = 549.75
$ws.Range("M86").Value = -265.7367999999999
$ws.Range("N86").Value = -2795.75

$ws.Range("H89").Value = 1242.826
$ws.Range("I89").Value = 1388.7368
$ws.Range("J89").Value = 549.75
$ws.Range("K89").Value = 6943.683999999999
$ws.Range("L89").Value = 2748.75
$ws.Range("M89").Value = -1327.683999999999
$ws.Range("N89").Value = -13980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2724.7932
$ws.Range("J16").Value = 2910.5
$ws.Range("L16").Value = 2910.5
$ws.Range("N16").Value = -3484.5

$ws.Range("H113").Value = 2724.7932
$ws.Range("J113").Value = 2910.5
$ws.Range("L113").Value = 2910.5
$ws.Range("N113").Value = -7250.5

$ws.Range("H122").Value = 5023.2354
$ws.Range("I122").Value = 4592.273
$ws.Range("J122").Value = 5813.3335
$ws.Range("K122").Value = 13776.819
$ws.Range("L122").Value = 17440.0005
$ws.Range("M122").Value = -11326.819
$ws.Range("N122").Value = -22340.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 7286.7856
$ws.Range("J23").Value = 16806.834
$ws.Range("L23").Value = 50420.50199999999
$ws.Range("N23").Value = -50890.50199999999

$ws.Range("H34").Value = 66984.47
$ws.Range("I34").Value = 131.66667
$ws.Range("J34").Value = 111553
$ws.Range("K34").Value = 395.00001
$ws.Range("L34").Value = 334659
$ws.Range("M34").Value = -311.00001
$ws.Range("N34").Value = -334827

$ws.Range("H37").Value = 148214.56
$ws.Range("J37").Value = 148214.56
$ws.Range("L37").Value = 444643.68
$ws.Range("N37").Value = -444867.68

$ws.Range("H92").Value = 186.16667
$ws.Range("I92").Value = 294.66666
$ws.Range("J92").Value = 150
$ws.Range("K92").Value = 883.9999799999999
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 364.0000200000001
$ws.Range("N92").Value = -2946

$ws.Range("H140").Value = 2788.2666
$ws.Range("I140").Value = 2702
$ws.Range("K140").Value = 8106
$ws.Range("M140").Value = -2926

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11057.462
$ws.Range("I70").Value = 7750
$ws.Range("J70").Value = 12527.444
$ws.Range("K70").Value = 7750
$ws.Range("L70").Value = 12527.444
$ws.Range("M70").Value = -7480
$ws.Range("N70").Value = -13067.444

$ws.Range("H73").Value = 11057.462
$ws.Range("I73").Value = 7750
$ws.Range("J73").Value = 12527.444
$ws.Range("K73").Value = 7750
$ws.Range("L73").Value = 12527.444
$ws.Range("M73").Value = -6814
$ws.Range("N73").Value = -14399.444

$ws.Range("H80").Value = 5210.7144
$ws.Range("I80").Value = 4429.1665
$ws.Range("J80").Value = 9900
$ws.Range("K80").Value = 4429.1665
$ws.Range("L80").Value = 9900
$ws.Range("M80").Value = -3431.1665
$ws.Range("N80").Value = -11896

$ws.Range("H83").Value = 5210.7144
$ws.Range("I83").Value = 4429.1665
$ws.Range("J83").Value = 9900
$ws.Range("K83").Value = 22145.8325
$ws.Range("L83").Value = 49500
$ws.Range("M83").Value = -17153.8325
$ws.Range("N83").Value = -59484

$ws.Range("H113").Value = 5722.222
$ws.Range("I113").Value = 4997.4
$ws.Range("J113").Value = 6628.25
$ws.Range("K113").Value = 4997.4
$ws.Range("L113").Value = 6628.25
$ws.Range("M113").Value = -2827.4
$ws.Range("N113").Value = -10968.25

$ws.Range("H132").Value = 7574.727
$ws.Range("I132").Value = 7065.375
$ws.Range("K132").Value = 21196.125
$ws.Range("M132").Value = -18666.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H16").Value = 1670.5
$ws.Range("I16").Value = 1670.5
$ws.Range("K16").Value = 1670.5
$ws.Range("M16").Value = -1500.5

$ws.Range("H82").Value = 11277.81
$ws.Range("I82").Value = 14486.167
$ws.Range("K82").Value = 14486.167
$ws.Range("M82").Value = -14125.167

$ws.Range("H85").Value = 11277.81
$ws.Range("I85").Value = 14486.167
$ws.Range("K85").Value = 14486.167
$ws.Range("M85").Value = -13238.167

$ws.Range("H93").Value = 16197
$ws.Range("I93").Value = 2645.9092
$ws.Range("J93").Value = 53462.5
$ws.Range("K93").Value = 2645.9092
$ws.Range("L93").Value = 53462.5
$ws.Range("M93").Value = -1397.9092
$ws.Range("N93").Value = -55958.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 29999
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H96").Value = 4049
$ws.Range("I96").Value = 3586
$ws.Range("J96").Value = 4975
$ws.Range("K96").Value = 3586
$ws.Range("L96").Value = 4975
$ws.Range("M96").Value = -2213
$ws.Range("N96").Value = -7721

$ws.Range("H107").Value = 982
$ws.Range("I107").Value = 982
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2946
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1026
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 4309.136
$ws.Range("I122").Value = 3250.1428
$ws.Range("J122").Value = 6162.375
$ws.Range("K122").Value = 9750.428400000001
$ws.Range("L122").Value = 18487.125
$ws.Range("M122").Value = -7300.428400000001
$ws.Range("N122").Value = -23387.125

$ws.Range("H129").Value = 72000
$ws.Range("J129").Value = 72000
$ws.Range("L129").Value = 72000
$ws.Range("N129").Value = -82000

$ws.Range("H132").Value = 4703.5
$ws.Range("I132").Value = 4750.4
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 14251.2
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -11721.2
$ws.Range("N132").Value = -17060
